# Add a "Save" column (H) to the s_vals sheet, matching the style of the
# existing header row, with a 0 value in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, bold/bordered/
# centered header style) onto the new H1 header cell, then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for row 2, plain (unstyled) numeric value like its neighbors.
$ws.Range("H2").Value = 0
